$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 3 data, with updated odds): Instituto vs Lanus @ 21:15:00
$ws.Range("C2").Value = "21:15:00"
$ws.Range("D2").Value = "Instituto"
$ws.Range("E2").Value = "Lanus"
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4.1
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 3.95
$ws.Range("J2").Value = 2.02
$ws.Range("K2").Value = 2.04
$ws.Range("P2").Value = 2.64
$ws.Range("Q2").Value = 1.59
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 4
$ws.Range("V2").Value = 1.37
$ws.Range("W2").Value = 1.31
$ws.Range("AC2").Value = 2.7
$ws.Range("AD2").Value = 4.9
$ws.Range("AE2").Value = 28
$ws.Range("AG2").Value = 6.2
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 60
$ws.Range("AL2").Value = 75
$ws.Range("AM2").Value = 290
$ws.Range("AN2").Value = 180
$ws.Range("AO2").Value = 130

# Row 3 (was row 4 data, with updated odds): Independiente Rivadavia vs Sarmiento de Junin @ 21:15:00
$ws.Range("D3").Value = "Independiente Rivadavia"
$ws.Range("E3").Value = "Sarmiento de Junin"
$ws.Range("F3").Value = 1.03
$ws.Range("G3").Value = 1.05
$ws.Range("H3").Value = 1.01
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 29
$ws.Range("K3").Value = 840
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.63
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 4.2
$ws.Range("T3").Value = 1.25
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 28
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 2.66
$ws.Range("AG3").Value = 7.6
$ws.Range("AH3").Value = 60
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 4.5
$ws.Range("AK3").Value = 13.5
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 10.5
$ws.Range("AO3").Value = 1000

# Row 4 (was row 5 data, with updated odds): Barranquilla vs Boyaca Patriotas @ 21:45:00, Colombian Primera B
$ws.Range("A4").Value = "Colombian Primera B"
$ws.Range("C4").Value = "21:45:00"
$ws.Range("D4").Value = "Barranquilla"
$ws.Range("E4").Value = "Boyaca Patriotas"
$ws.Range("F4").Value = 1.76
$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 6.6
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.55
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 3.9
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 1.63
$ws.Range("Q4").Value = 2.48
$ws.Range("R4").Value = 1.16
$ws.Range("S4").Value = 6.6
$ws.Range("T4").Value = 1.66
$ws.Range("U4").Value = 2.36
$ws.Range("V4").Value = 1.16
$ws.Range("W4").Value = 2.24
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 4.2
$ws.Range("AC4").Value = 4.5
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 55
$ws.Range("AF4").Value = 8.6
$ws.Range("AG4").Value = 9.2
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 110
$ws.Range("AJ4").Value = 44
$ws.Range("AK4").Value = 55
$ws.Range("AL4").Value = 130
$ws.Range("AM4").Value = 470
$ws.Range("AN4").Value = 110
$ws.Range("AO4").Value = 180

# Remove old row 5 entirely (shift rows up)
$ws.Rows("5:5").Delete()
